$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Signups": the hidden "Deleted_ý" flag column (A) and the "Id" column
# (B) switch from boolean/numeric storage to plain text, and the
# "Party size" column (E) switches from numeric storage to text. The
# leftover per-cell number-format override inherited from column B's
# (text) column style is cleared back to the default ("Normal") style so
# the cell again relies on the column-level style.
# ---------------------------------------------------------------------------
$wsSignups = $wb.Worksheets.Item("Signups")

function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Row 2 - Brice
Set-TextValue $wsSignups.Range("A2") "0"
Set-TextValue $wsSignups.Range("B2") "1"
Set-TextValue $wsSignups.Range("E2") "3"

# Row 3 - Ryan
Set-TextValue $wsSignups.Range("A3") "0"
Set-TextValue $wsSignups.Range("B3") "2"
Set-TextValue $wsSignups.Range("E3") "4"

# Row 4 - David
Set-TextValue $wsSignups.Range("A4") "0"
Set-TextValue $wsSignups.Range("B4") "3"
Set-TextValue $wsSignups.Range("E4") "2"

# ---------------------------------------------------------------------------
# Sheet "Log": the prior "Updated signup" audit entry is replaced by the
# entry that originally recorded Brice's signup, and three further audit
# rows are appended recording the Ryan signup, the follow-up phone-number
# correction for Ryan, and the David signup.
# ---------------------------------------------------------------------------
$wsLog = $wb.Worksheets.Item("Log")

$wsLog.Range("A2").Value = "ea3bdcd5"
$wsLog.Range("B2").Value = "John.Doe"
$wsLog.Range("C2").Value = 45888.5568055556
Set-TextValue $wsLog.Range("D2") "1"
$wsLog.Range("E2").Value = "Added signup: Name: Brice, Phone: 555-5551, Party Size: 3"

Set-TextValue $wsLog.Range("A3") "33427300"
$wsLog.Range("B3").Value = "John.Doe"
$wsLog.Range("C3").Value = 45888.559282407397
Set-TextValue $wsLog.Range("D3") "2"
$wsLog.Range("E3").Value = "Added signup: Name: Ryan, Phone: 555-5559, Party Size: 4"

$wsLog.Range("A4").Value = "8187ebee"
$wsLog.Range("B4").Value = "John.Doe"
$wsLog.Range("C4").Value = 45888.480891203704
Set-TextValue $wsLog.Range("D4") "2"
$wsLog.Range("E4").Value = "Updated signup: Name: [unchanged], Phone: 555-5552, Party Size: [unchanged]"

Set-TextValue $wsLog.Range("A5") "25548527"
$wsLog.Range("B5").Value = "John.Doe"
$wsLog.Range("C5").Value = 45888.482013888897
Set-TextValue $wsLog.Range("D5") "3"
$wsLog.Range("E5").Value = "Added signup: Name: David, Phone: 555-5553, Party Size: 2"

# The Timestamp column's display format drops seconds and the 4-digit year
# once the log is populated with real data.
$wsLog.Range("C2:C5").NumberFormat = "d/mm/yy\ h:mm;@"
